# Applies the "copy old investment folder" update to the 总结 (Summary) sheet.
# 基金 / 理财 sheets are untouched (their only diffs are cached <v> values for
# the volatile NOW() formula, which recalculates automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("总结")

# ---- Row 6 ----------------------------------------------------------------
$ws.Range("M6").Value = 8443.75
$ws.Range("O6").Formula = "=SUM(L:L) +M6"

$ws.Range("T6").Value = 6079.31
$ws.Range("U6").Value = 9596.65
# W6 / AD6 lose their SUM() formulas and become plain numbers.
$ws.Range("W6").Value = 19046.560000000001
$ws.Range("AA6").Value = 735.78
$ws.Range("AD6").Value = 2403.2399999999998

# ---- Row 7 ------------------------------------------------------------
# O7 becomes a formula (was a literal number).
$ws.Range("O7").Formula = "=813077.18-143160.54"

# The old T7 formula (-22.47*2) shifts left into S7; T7 gets a new plain value.
$ws.Range("S7").Formula = "=-22.47*2"
$ws.Range("T7").Value = 2750.68
# U7:U8 share one relative formula "=T7+S7".
$ws.Range("U7:U8").Formula = "=T7+S7"

$ws.Range("W7").Value = 201676.06
$ws.Range("Z7").Value = 141000
$ws.Range("AA7").Value = 480.06
$ws.Range("AD7").Value = 311000

# ---- Row 8 ------------------------------------------------------------
# The old T8 formula (-31.45) shifts left into S8; T8 gets a new plain value.
$ws.Range("S8").Formula = "=-31.45"
$ws.Range("T8").Value = 2054.7199999999998

$ws.Range("AA8").Value = 651.66999999999996

# ---- Row 9 (new entries in the Q:U block) ---------------------------------
$ws.Range("Q9").Value = "交银新生活力灵活配置混合"
$ws.Range("R9").Value = 27266.98
$ws.Range("S9").Formula = "=-18.37-30"
$ws.Range("T9").Value = 2.6
$ws.Range("U9").Formula = "=T9+S9"

# ---- Row 10 -----------------------------------------------------------
$ws.Range("Q10").Value = "中银珍利混合C"
$ws.Range("R10").Value = 18879.490000000002
$ws.Range("S10").Value = 0
# U10:U12 share one relative formula "=T10+S10".
$ws.Range("U10:U12").Formula = "=T10+S10"

# ---- Row 11 -----------------------------------------------------------
$ws.Range("Q11").Value = "长信乐信灵活配置混合C"
$ws.Range("R11").Value = 11326.75
$ws.Range("S11").Value = 0

# ---- Row 12 -----------------------------------------------------------
$ws.Range("Q12").Value = "广发趋势优选灵活配置混合A"
$ws.Range("R12").Value = 7544.21
$ws.Range("S12").Formula = "=-R12*0.15/100"
$ws.Range("S12").NumberFormat = "0.00"

# ---- New rows 19-21 (additional 聚益生金 purchases) ------------------------
$ws.Range("A19").Value = "聚益生金63天C款"
$ws.Range("B19").Value = 43971
$ws.Range("B19").NumberFormat = "m/d/yy"
$ws.Range("C19").Value = 50000
$ws.Range("D19").Value = 345.2

$ws.Range("A20").Value = "聚益生金91天A款"
$ws.Range("B20").Value = 43976
$ws.Range("B20").NumberFormat = "m/d/yy"
$ws.Range("C20").Value = 90000
$ws.Range("D20").Value = 845.55

$ws.Range("A21").Value = "聚益生金91天A款"
$ws.Range("B21").Value = 44001
$ws.Range("B21").NumberFormat = "m/d/yy"
$ws.Range("C21").Value = 160000
$ws.Range("D21").Value = 1049.44

# ---- Updated totals in rows 26-28 ------------------------------------------
$ws.Range("D26").Value = 6889.49
$ws.Range("D27").Value = 876.18
$ws.Range("D28").Value = 1492.84

# ---- Column widths ----------------------------------------------------
# Column Q (17) widens; column S (19) gains an explicit (best-fit) width.
$ws.Columns.Item(17).ColumnWidth = 24
$ws.Columns.Item(19).ColumnWidth = 11.833333333333332

# ---- Selection moves to O8 -------------------------------------------------
$ws.Range("O8").Select()
